$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (header "Förändrad") holds a serial date value that was bumped
# by one day (45177 -> 45178) for every data row (rows 2 through 171).
for ($row = 2; $row -le 171; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45177) {
        $cell.Value2 = 45178
    }
}
